$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("2010")
$ws = $wb.Worksheets.Item("2010-18")

# Make the "2010-18" sheet the active/selected tab (matches workbook.xml
# activeTab and sheet1/sheet3 tabSelected flip in the diff).
$ws.Activate()

# Insert a new column before column I ("AET") to hold the new
# " water added by FlowModel (mm)" series, shifting existing data right.
$ws.Columns.Item(9).Insert()

# New column header text (row 1).
$ws.Range("I1").Value = " water added by FlowModel (mm)"

# New calibration run data (row 4): David Holman's latest PEST
# calibration of Mehama37 - "CW3M ~C401" / "Demo_Baseline 5/20/21".
$ws.Range("B4").Value = "Demo_Baseline 5/20/21"
$ws.Range("A4").Value = "CW3M ~C401"
$ws.Range("C4").Value = "2010-18"

$ws.Range("D4").Value = 1150.0584241111112
$ws.Range("E4").Value = 1612.6987305555554
$ws.Range("F4").Value = 14.557748111111115
$ws.Range("G4").Value = 52.671807666666659
$ws.Range("H4").Value = 5.2565644444444439
$ws.Range("I4").Value = 8.9746754444444452
$ws.Range("J4").Value = 2.782013222222222
$ws.Range("K4").Value = 609.74378122222208
$ws.Range("L4").Value = 44.391417555555549
$ws.Range("M4").Value = 1018.5836656666668
$ws.Range("N4").Value = 1167.9530299999999
$ws.Range("O4").Value = 517762.13888888888
$ws.Range("P4").Value = 286785.73958333331
$ws.Range("Q4").Value = -0.76404255555555545
$ws.Range("R4").Value = -0.00027444444444444445

# Match the formatting used by the surrounding table: copy number-format
# styles from analogous cells (within this sheet where possible, and from
# the "2010" sheet for the two styles -- index 6 and 7 -- that aren't used
# anywhere on "2010-18" yet).
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)

$ws1.Range("E4").Copy()
$ws.Range("D4:G4").PasteSpecial(-4122)
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("M4:N4").PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("Q4").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("J4").PasteSpecial(-4122)

$ws1.Range("N8").Copy()
$ws.Range("O4").PasteSpecial(-4122)

$ws.Range("O3").Copy()
$ws.Range("P4").PasteSpecial(-4122)

$ws.Range("R2").Copy()
$ws.Range("R4").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# New column's width (col R, the rightmost "mass balance discrepancy
# (fraction)" column after the insert) gets its own explicit width.
$ws.Columns.Item(18).ColumnWidth = 10.5546875

# Final selection left on the sheet after the edit.
$ws.Range("A6:XFD6").Select()
